$d = $word.ActiveDocument

function Replace-Span([string]$old, [string]$new) {
    $found = $d.Content.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Host "WARNING: span not found: [$old]"
    }
    return $found
}

# 1) Title paragraph: " - Ron Peer and Lavi Batzia " + "- " + " ""  + "mapPIN" + """ -> merge, drop proofErr
$old1 = " – Ron Peer and Lavi Batzia –  “mapPIN”"
Replace-Span $old1 $old1

# 2) "How to demo" (Transfer Turns Between players): drop "the and" gramStart proofErr, merge runs
$old2 = ": Once a player has made their choice and the score has been calculated, it is shown on the and screen and then"
Replace-Span $old2 $old2

# 3) "How to demo" (Distance calculation - city mode): drop "interval from the timer" gramStart proofErr, merge runs
$old3 = ", the coordinates of the mouse along with the time interval from the timer will be used to calculate distance between the mouse and city of the round"
Replace-Span $old3 $old3

# 4) "How to demo" (Game over screen): drop "a" spellStart/gramStart proofErr, merge runs
$old4 = " containing a"
Replace-Span $old4 $old4

# 5) "How to demo" (player object): drop "matter" gramStart proofErr, merge runs
$old5 = "During each of the players turns in the game, no matter the mode,"
Replace-Span $old5 $old5

# 6) Version control paragraph: drop multiple proofErr markers, merge runs
$old6 = "basic and all necessary commands in order to be able to collaboratively work on the same project and track changes. Using branching, each one of us was able to focus on their on part while still updating the other, using push-pull commands to the origin. During the learning process, we explored the previous commits we made due to changes and bugs. We merged into the main branch a few times when full tested code was finished. As part of the development, we frequently switched branches to learn, test, improve and integrate each other’s work."
Replace-Span $old6 $old6

# 7) "git add" bullet: drop gramStart proofErr, merge runs
$old7 = "git add"
Replace-Span $old7 $old7

# 8) Insert a new bullet "git log (q to exit)" right after "git stash" (and before "git stash pop")
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text.TrimEnd([char]13, [char]7) -eq "git stash") {
        $target = $para
        break
    }
}
if ($target -ne $null) {
    $target.Range.InsertParagraphAfter()
    $newPara = $target.Next()
    $newPara.Range.Text = "git log (q to exit)"
} else {
    Write-Host "WARNING: 'git stash' paragraph not found"
}

Write-Host "Done"
